$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00000000181285514432328
$ws.Range("C2").Value = 0.000000001184464077636434
$ws.Range("D2").Value = 0.000000001864478199970915
$ws.Range("E2").Value = 0.9999999951382343
$ws.Range("F2").Value = 46080
